# Updates to get dispatch and vehicle shares working
# Sets the BAU Guaranteed Dispatch Percentage (column B) to 1 (100%) for:
#   onshore wind, solar PV, solar thermal, biomass, offshore wind
# The remaining year columns (C:AK) contain formulas referencing column B
# of the same row, so they recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

$ws.Range("B6").Value = 1   # onshore wind
$ws.Range("B7").Value = 1   # solar PV
$ws.Range("B8").Value = 1   # solar thermal
$ws.Range("B9").Value = 1   # biomass
$ws.Range("B14").Value = 1  # offshore wind

# Make BGDPbES the active/selected sheet, matching the workbook view update
$ws.Activate()
$ws.Range("B17").Select()
